$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = "Token has expired"
$ws.Cells.Item(38, 3).Value = 1
$ws.Range("I26").Select()
